$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# --- Add new TODO rows (132-140) ---------------------------------------
# Row 132 stays blank (a spacer row between the existing list and the new
# "Post SAM 2014.11.24 Release" section).

# Row 133: section header
$ws.Range("A133").Value = "Post SAM 2014.11.24 Release"

# Rows 134-140: new post-release bugs/issues
$ws.Range("A134").Value = "Not Done"
$ws.Range("B134").Value = "Display issues on very high resolution screens (see screenshots from Ted James)"
$ws.Range("C134").Value = "Aron"

$ws.Range("A135").Value = "Not Done"
$ws.Range("B135").Value = "Change curtailment and availability factors to percentages to be consistent with other losses"
$ws.Range("C135").Value = "Janine"

$ws.Range("A136").Value = "Not Done"
$ws.Range("B136").Value = "Mystery sam.exe crash on some Windows computers"
$ws.Range("C136").Value = "Aron"

$ws.Range("A137").Value = "Not Done"
$ws.Range("B137").Value = "Remove lk autosave to avoid junk temporary files"
$ws.Range("C137").Value = "Aron"

$ws.Range("A138").Value = "Not Done"
$ws.Range("B138").Value = "Add monthly scaling to input hourly load page (monthly scaling available for residential belpe but not commercial)"
$ws.Range("C138").Value = "Janine"
$ws.Rows.Item(138).RowHeight = 15

$ws.Range("A139").Value = "Not Done"
$ws.Range("B139").Value = "Physical trough default T_startup (lowered from 300 to 250 to avoid simulation warnings, but causes 3% annual output reduction)"
$ws.Range("C139").Value = "Steve/Ty"
$ws.Rows.Item(139).RowHeight = 30

$ws.Range("A140").Value = "Not Done"
$ws.Range("B140").Value = "Add LHS and stepwise script functions to LK"
$ws.Range("C140").Value = "Aron"

# --- Formatting for the new rows ----------------------------------------
# Column B wraps text; column A keeps the existing centered style already
# used throughout the sheet (applied automatically via the column style).
$ws.Range("B134:B140").WrapText = $true

# Header row (133) - bold white text on accent fill for A133, bold black
# text on the same accent fill for B133:E133.
$headerRange = $ws.Range("A133:E133")
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = -0.249977111117893
$ws.Range("A133").Font.Bold = $true
$ws.Range("A133").Font.ThemeColor = 2
$ws.Range("A133").HorizontalAlignment = -4131
$ws.Range("B133:E133").Font.Bold = $true

# --- Move the active view from "Project Ideas" to "To Do" ---------------
$ws.Range("A141").Select()
$ws.Activate()
